$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 291
$ws1.Range("F4").Value = 2394
$ws1.Range("F5").Value = 1771
$ws1.Range("F6").Value = 340
$ws1.Range("F7").Value = 98
$ws1.Range("F8").Value = 827
$ws1.Range("F9").Value = 166

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 291
$ws4.Range("F4").Value = 2394
$ws4.Range("F5").Value = 1771
$ws4.Range("F6").Value = 340
$ws4.Range("F8").Value = 98
$ws4.Range("F9").Value = 827
$ws4.Range("F10").Value = 166
